$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns for rows 2-51 with the latest
# scrape data. Column D values are price strings that sometimes look numeric
# (e.g. "1.004"); to keep them stored as text (matching the original cells) we
# briefly switch the cell to Text format before assigning the value, then restore
# the default "Normal" style so no formatting differences are introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.933.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5083"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.643.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5463"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7867"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.993.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.437"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.974"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.047"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.886"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.891"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.239"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05040"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.195"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.544"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.371"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.606"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.135.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.007"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.548"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.640"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8169"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.779.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4526"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05075"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "

# Rows 39 and 40 swapped order (VeChain <-> BabyDogeCoin) with refreshed values
$ws.Range("B39").Value = "BabyDogeCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₈134"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +15.18%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
